$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    [double]"0.9999883732771242",
    [double]"0.9990763349190175",
    [double]"0.9999999999998855",
    [double]"0.9999989859190584",
    [double]"0.9999994160120935",
    [double]"1.085303545036088e-05",
    [double]"0.000862200808886027",
    [double]"8.958573410540517e-14",
    [double]"1.082021660499151e-06",
    [double]"5.410108750424426e-07",
    [double]"0.0001805656912053873",
    [double]"0.00329439454989242",
    [double]"0.9999069862169936",
    [double]"0.003434643746218193",
    [double]"64.86213150418904",
    [double]"90.45852382642124"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
